$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new data rows before the current last data row (19) ---
# Excel shifts the footer rows (24/25 -> 26/27) and merged cells automatically.
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(19).Insert()

# Copy the formatting (borders/fill/font/number format) of row 18 onto the
# two freshly inserted rows 19 and 20 so they match the other "middle" data
# rows instead of the generic default style Insert() leaves behind.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J20").PasteSpecial(-4122)

# --- Update header / summary figures ---
$ws.Range("E11").Value = 254332
$ws.Range("F13").Value = 3

# --- Rewrite the worker/period detail table (rows 16-21) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1007126667"
$ws.Range("D16").Value = "LUIS DANIEL DE LA ROSA PEREZ"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 13286
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235039810"
$ws.Range("D17").Value = "JESUS ALBERTO ESCORCIA SANTIAGO"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 13286
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007126667"
$ws.Range("D18").Value = "LUIS DANIEL DE LA ROSA PEREZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1235039810"
$ws.Range("D19").Value = "JESUS ALBERTO ESCORCIA SANTIAGO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1007126667"
$ws.Range("D20").Value = "LUIS DANIEL DE LA ROSA PEREZ"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1235039810"
$ws.Range("D21").Value = "JESUS ALBERTO ESCORCIA SANTIAGO"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500
